# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates column G ("K") values for rows 2-22 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3
    3  = 8
    4  = 5
    5  = 2
    6  = 6
    7  = 5
    8  = 1
    9  = 6
    10 = 7
    11 = 3
    12 = 2
    13 = 2
    14 = 4
    15 = 8
    16 = 4
    17 = 3
    18 = 3
    19 = 4
    20 = 6
    21 = 2
    22 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
